$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("F5").Value = 2.8
$ws.Range("G5").Value = 3.5
$ws.Range("H5").Value = 2.32
$ws.Range("I5").Value = 2.8
$ws.Range("Q5").Value = 1.73

# Row 7 updates
$ws.Range("F7").Value = 1.43
$ws.Range("H7").Value = 9.199999999999999
$ws.Range("I7").Value = 12
$ws.Range("N7").Value = 2.96
$ws.Range("T7").Value = 2.58
$ws.Range("U7").Value = 1.55
$ws.Range("V7").Value = 1.09
$ws.Range("Z7").Value = 120
$ws.Range("AI7").Value = 270

# Row 11 updates
$ws.Range("G11").Value = 2.48

# Row 13 updates
$ws.Range("P13").Value = 1.99
$ws.Range("Q13").Value = 1.82
